# Append newly tracked airline flight data for 2021-01-02 through 2021-01-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-01-02", 57, 54),
    @("2021-01-03", 52, 50),
    @("2021-01-04", 55, 49),
    @("2021-01-05", 45, 43),
    @("2021-01-06", 41, 40),
    @("2021-01-07", 53, 50),
    @("2021-01-08", 48, 47),
    @("2021-01-09", 42, 41),
    @("2021-01-10", 48, 46),
    @("2021-01-11", 39, 37),
    @("2021-01-12", 49, 47),
    @("2021-01-13", 47, 45),
    @("2021-01-14", 48, 48),
    @("2021-01-15", 38, 38),
    @("2021-01-16", 47, 46),
    @("2021-01-17", 36, 36),
    @("2021-01-18", 44, 43),
    @("2021-01-19", 41, 37)
)

$startRow = 272
$lastRow = $startRow + $data.Count - 1

# Carry the existing row's number/text formatting down into the new rows
# (same as dragging the fill handle from the last populated row).
$ws.Range("A271:D271").Copy() | Out-Null
$ws.Range("A" + $startRow + ":D" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Formula = "=C" + $row + "/B" + $row
}

$ws.Range("I279").Select()
